# Generate Report for Handback
# Adds two new handed-back files (e43d5d25-... and eeb270db-...) as new rows
# to the "Overview", "zh-cn" and "de-de" worksheets of the handback-status
# workbook, mirroring the existing rows' structure (values + hyperlinks).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared data for the two newly handed-back files
# ---------------------------------------------------------------------
$status = "Handed back: in sync with en-US"

$file1 = "e43d5d25-df8c-4566-982f-7558f8ff7fed"
$file2 = "eeb270db-1116-4259-8295-a1fef7cd0f3c"

$hash1 = "ed6237191c14abd94fde32494a13d9f14b3afb24"
$hash2 = "f2a5ae88b4f1432276130e842422c4709a319429"

$md1 = "$file1.md"
$md2 = "$file2.md"

$xlfZh1 = "$file1.$hash1.zh-cn.xlf"
$xlfZh2 = "$file2.$hash2.zh-cn.xlf"
$xlfDe1 = "$file1.$hash1.de-de.xlf"
$xlfDe2 = "$file2.$hash2.de-de.xlf"

# source-repo commit used to build the "Col A / Col E" https://.../e2e/*.md links
$srcCommit1 = "956a4ee626bf11ce52af36561090510fc7d677f4"
$srcCommit2 = "ae699ac8d38d6d7e3b6b93e5894e5fc3a8079193"

$zhSourceCommit1 = "956a4ee626bf11ce52af36561090510fc7d677f4"
$zhSourceCommit2 = "ae699ac8d38d6d7e3b6b93e5894e5fc3a8079193"
$deSourceCommit1 = "956a4ee626bf11ce52af36561090510fc7d677f4"
$deSourceCommit2 = "ae699ac8d38d6d7e3b6b93e5894e5fc3a8079193"

$zhHandoffCommit1 = "aa392ea731d5c4e557c91b8136875b286450876b"
$zhHandoffCommit2 = "7172e316e8456f90d0522a5ce4641c778a231498"
$deHandoffCommit1 = "aa392ea731d5c4e557c91b8136875b286450876b"
$deHandoffCommit2 = "7172e316e8456f90d0522a5ce4641c778a231498"

$zhTargetCommit1 = "910385770b562b5285de4d5e4943c4c1e68c3e85"
$zhTargetCommit2 = "7010e27fe733a7d8dbe55aa5bf7b8d316c01f42a"
$deTargetCommit1 = "910385770b562b5285de4d5e4943c4c1e68c3e85"
$deTargetCommit2 = "7010e27fe733a7d8dbe55aa5bf7b8d316c01f42a"

$zhHandbackCommit1 = "4ac43d427d32fe52b4553f413aafb3ce1691999d"
$zhHandbackCommit2 = "73c982bcb1943c0eb5c31f9bd5e4baccebf3912a"
$deHandbackCommit1 = "4ac43d427d32fe52b4553f413aafb3ce1691999d"
$deHandbackCommit2 = "73c982bcb1943c0eb5c31f9bd5e4baccebf3912a"

# Timestamps (zh-cn)
$zhHandoffTime1 = "2016-02-18 10:26:19"
$zhHandbackTime1 = "2016-02-18 10:27:15"
$zhHandoffTime2 = "2016-02-18 10:26:19"
$zhHandbackTime2 = "2016-02-18 10:27:15"

# Timestamps (de-de)
$deHandoffTime1 = "2016-02-18 10:26:35"
$deHandbackTime1 = "2016-02-18 10:27:38"
$deHandoffTime2 = "2016-02-18 10:26:35"
$deHandbackTime2 = "2016-02-18 10:27:38"

$reason = "Include"

# =======================================================================
# Sheet "Overview": File Name | zh-cn | de-de
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Hyperlinks.Add($ov.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit1/e2e/$md1", "", "", $md1)
$ov.Range("B6").Value2 = $status
$ov.Range("C6").Value2 = $status

$ov.Hyperlinks.Add($ov.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit2/e2e/$md2", "", "", $md2)
$ov.Range("B7").Value2 = $status
$ov.Range("C7").Value2 = $status

# =======================================================================
# Sheet "zh-cn": Source File Name | Status | Correspond Handoff File |
#   Correspond Handoff Datetime | Target File | Correspond Handback File |
#   Correspond Handback DateTime | Handoff Reason | Dependency From
# =======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# Row 6 -> file1
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$zhSourceCommit1/e2e/$md1", "", "", $md1)
$zh.Range("B6").Value2 = $status
$zh.Hyperlinks.Add($zh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffCommit1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh1", "", "", $xlfZh1)
$zh.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("D6").Value2 = $zhHandoffTime1
$zh.Hyperlinks.Add($zh.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$zhTargetCommit1/e2e/$md1", "", "", $md1)
$zh.Hyperlinks.Add($zh.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$zhHandbackCommit1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh1", "", "", $xlfZh1)
$zh.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G6").Value2 = $zhHandbackTime1
$zh.Range("H6").Value2 = $reason

# Row 7 -> file2
$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$zhSourceCommit2/e2e/$md2", "", "", $md2)
$zh.Range("B7").Value2 = $status
$zh.Hyperlinks.Add($zh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffCommit2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh2", "", "", $xlfZh2)
$zh.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("D7").Value2 = $zhHandoffTime1
$zh.Hyperlinks.Add($zh.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$zhTargetCommit2/e2e/$md2", "", "", $md2)
$zh.Hyperlinks.Add($zh.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$zhHandbackCommit2/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$xlfZh2", "", "", $xlfZh2)
$zh.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("G7").Value2 = $zhHandbackTime1
$zh.Range("H7").Value2 = $reason

# =======================================================================
# Sheet "de-de": same columns as "zh-cn"
# =======================================================================
$de = $wb.Worksheets.Item("de-de")

# Row 6 -> file1
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/$deSourceCommit1/e2e/$md1", "", "", $md1)
$de.Range("B6").Value2 = $status
$de.Hyperlinks.Add($de.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffCommit1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe1", "", "", $xlfDe1)
$de.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("D6").Value2 = $deHandoffTime1
$de.Hyperlinks.Add($de.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$deTargetCommit1/e2e/$md1", "", "", $md1)
$de.Hyperlinks.Add($de.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$deHandbackCommit1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe1", "", "", $xlfDe1)
$de.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G6").Value2 = $deHandbackTime1
$de.Range("H6").Value2 = $reason

# Row 7 -> file2
$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/$deSourceCommit2/e2e/$md2", "", "", $md2)
$de.Range("B7").Value2 = $status
$de.Hyperlinks.Add($de.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffCommit2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe2", "", "", $xlfDe2)
$de.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("D7").Value2 = $deHandoffTime1
$de.Hyperlinks.Add($de.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$deTargetCommit2/e2e/$md2", "", "", $md2)
$de.Hyperlinks.Add($de.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$deHandbackCommit2/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$xlfDe2", "", "", $xlfDe2)
$de.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("G7").Value2 = $deHandbackTime1
$de.Range("H7").Value2 = $reason
